$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 'Menu', 15, 5, 20),
    @(3, 'Nav', 9, 1, 10),
    @(4, 'DataGrid', 4, 6, 10),
    @(5, 'Tree', 9, 0, 9),
    @(6, 'Popover', 6, 2, 8),
    @(7, 'Dialog', 7, 1, 8),
    @(8, 'Combobox', 3, 3, 6),
    @(9, 'Table', 5, 1, 6),
    @(10, 'TagPicker', 4, 1, 5),
    @(11, 'Dropdown', 2, 2, 4),
    @(12, 'Tooltip', 2, 2, 4),
    @(13, 'Toolbar', 3, 1, 4),
    @(14, 'Virtualizer', 4, 0, 4),
    @(15, 'Calendar Compat', 3, 0, 3),
    @(16, 'Skeleton', 2, 1, 3),
    @(17, 'MessageBar', 3, 0, 3),
    @(18, 'TeachingPopover', 3, 0, 3),
    @(19, 'Drawer', 2, 0, 2),
    @(20, 'DatePicker', 0, 2, 2),
    @(21, 'DatePickerCompat', 1, 1, 2),
    @(22, 'Accordion', 2, 0, 2),
    @(23, 'Switch', 2, 0, 2),
    @(24, 'List', 2, 0, 2),
    @(25, 'Toast', 1, 1, 2),
    @(26, 'Portal', 2, 0, 2),
    @(27, 'Tabs', 2, 0, 2),
    @(28, 'TabList', 1, 1, 2),
    @(29, 'Slider', 0, 2, 2),
    @(30, 'FluentProvider', 0, 2, 2),
    @(31, 'Input', 0, 1, 1),
    @(32, 'Avatar', 1, 0, 1),
    @(33, 'Checkbox', 1, 0, 1),
    @(34, 'Badge', 1, 0, 1),
    @(35, 'AvatarGroup', 1, 0, 1),
    @(36, 'FocusTrapZone', 1, 0, 1),
    @(37, 'Popup', 1, 0, 1),
    @(38, 'SwatchPicker', 0, 1, 1),
    @(39, 'Image', 0, 1, 1),
    @(40, 'MenuItem', 1, 0, 1),
    @(41, 'Label', 1, 0, 1),
    @(42, 'Textarea', 1, 0, 1),
    @(43, 'Button', 0, 1, 1),
    @(44, 'Keytip', 0, 0, 0),
    @(45, 'Segment', 0, 0, 0),
    @(46, 'InfoLabel', 0, 0, 0),
    @(47, 'Tag', 0, 0, 0),
    @(48, 'SplitButton', 0, 0, 0),
    @(49, 'ColorPicker', 0, 0, 0),
    @(50, 'SpinButton', 0, 0, 0),
    @(51, 'Rating', 0, 0, 0),
    @(52, 'Pickers', 0, 0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

Write-Output "Updated component stats table"
